$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-02-22 14:25:01"
$wsZh.Range("G5").Value = "2016-02-22 14:25:48"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-02-22 14:25:15"
$wsDe.Range("G5").Value = "2016-02-22 14:26:10"
